# GED_Python_Workgroup_2018_01_23.pptx — "updated slides with agenda"
#
#   1. The footer "date" placeholder (the datetimeFigureOut field showing
#      1/4/2018) is refreshed to 1/11/2018 on the slide master and on every
#      slide layout.
#   2. The first agenda bullet on the "Workgroup 3" slide gets "(finish)"
#      appended: "Lesson 5 - Lists and dictionaries" -> "... (finish)".

$p = $ppt.ActivePresentation

$ppPlaceholderDate = 16
$oldDate = "1/4/2018"
$newDate = "1/11/2018"

function Update-DatePlaceholder {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePh = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePh = $true
            }
        } catch {}
        if (-not $isDatePh -and $sh.Name -like "Date Placeholder*") {
            $isDatePh = $true
        }

        if ($isDatePh -and $sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# --- 1. Refresh the date placeholder on the slide master ... ---
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# --- ... and on every slide layout ---
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# --- 2. "Lesson 5" now reads "... (finish)" on the Workgroup 3 slide ---
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            $bodyRange = $shape.TextFrame.TextRange
            $firstPara = $bodyRange.Paragraphs(1, 1)
            $firstParaText = $firstPara.Text.TrimEnd("`r")
            if ($firstParaText -eq "Lesson 5 – Lists and dictionaries") {
                $run = $bodyRange.Characters($firstPara.Start, $firstPara.Length)
                $run.Text = "Lesson 5 – Lists and dictionaries (finish)"
            }
        }
    }
}

Write-Output "done"
